# Using Yamakawa's method instead of the 2 FAMM approach
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tweak a handful of existing data points in the X FAMM (bottom) table ---
$ws.Range("C14").Value = -0.7
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0.4

$ws.Range("B17").Value = -0.4
$ws.Range("C17").Value = 0
$ws.Range("E17").Value = 0.7

# --- New block: row 25 section titles ---
$ws.Range("C25").Value = "YAMAKAWAWAWA"
$ws.Range("K25").Value = "YAMAKAWAWAWAW FAMM"
$ws.Range("N25").Value = "POSITION"

# --- Row 26: column headers for the small X table + the new Yamakawa FAMM table ---
$ws.Range("A26").Value = "X"
$ws.Range("B26").Value = "nl"
$ws.Range("C26").Value = "ns"
$ws.Range("D26").Value = "ze"
$ws.Range("E26").Value = "ps"
$ws.Range("F26").Value = "pl"

$ws.Range("K26").Value = "FAMM3"
$ws.Range("L26").Value = "NL"
$ws.Range("M26").Value = "NS"
$ws.Range("N26").Value = "ZE"
$ws.Range("O26").Value = "PS"
$ws.Range("P26").Value = "PL"
$ws.Range("K26:P26").Font.Bold = $true
$ws.Range("K26:P26").HorizontalAlignment = -4108

# --- Row 27 ---
$ws.Range("A27").Value = "a"
$ws.Range("B27").Value = -4
$ws.Range("C27").Value = -3
$ws.Range("D27").Value = -1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 0

$ws.Range("K27").Value = "NL"
$ws.Range("K27").Font.Bold = $true
$ws.Range("K27").HorizontalAlignment = -4108
$ws.Range("L27").Value = "NVL"
$ws.Range("M27").Value = "NL"
$ws.Range("N27").Value = "NL"
$ws.Range("O27").Value = "NS"
$ws.Range("P27").Value = "NS"

# --- Row 28 ---
$ws.Range("A28").Value = "b"
$ws.Range("B28").Value = -4
$ws.Range("C28").Value = -2
$ws.Range("D28").Value = -0.5
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 1

$ws.Range("K28").Value = "NS"
$ws.Range("K28").Font.Bold = $true
$ws.Range("K28").HorizontalAlignment = -4108
$ws.Range("L28").Value = "NL"
$ws.Range("M28").Value = "NM"
$ws.Range("N28").Value = "NS"
$ws.Range("O28").Value = "NS"
$ws.Range("P28").Value = "NS"

# --- Row 29 ---
$ws.Range("A29").Value = "c"
$ws.Range("B29").Value = -3
$ws.Range("C29").Value = -1
$ws.Range("D29").Value = 0.5
$ws.Range("E29").Value = 2
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 1

$ws.Range("J29").Value = "ANGLE"

$ws.Range("K29").Value = "ZE"
$ws.Range("K29").Font.Bold = $true
$ws.Range("K29").HorizontalAlignment = -4108
$ws.Range("L29").Value = "NS"
$ws.Range("M29").Value = "NS"
$ws.Range("N29").Value = "ZE"
$ws.Range("O29").Value = "PS"
$ws.Range("P29").Value = "PS"

# --- Row 30 ---
$ws.Range("A30").Value = "d"
$ws.Range("B30").Value = -2
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 3
$ws.Range("F30").Value = 4
$ws.Range("G30").Value = 0

$ws.Range("K30").Value = "PS"
$ws.Range("K30").Font.Bold = $true
$ws.Range("K30").HorizontalAlignment = -4108
$ws.Range("L30").Value = "PS"
$ws.Range("M30").Value = "PS"
$ws.Range("N30").Value = "PS"
$ws.Range("O30").Value = "PM"
$ws.Range("P30").Value = "PL"

# --- Row 31 ---
$ws.Range("K31").Value = "PL"
$ws.Range("K31").Font.Bold = $true
$ws.Range("K31").HorizontalAlignment = -4108
$ws.Range("L31").Value = "PS"
$ws.Range("M31").Value = "PS"
$ws.Range("N31").Value = "PL"
$ws.Range("O31").Value = "PL"
$ws.Range("P31").Value = "PVL"

# --- Row 32: headers for the new small Y table ---
$ws.Range("A32").Value = "Y"
$ws.Range("B32").Value = "nl"
$ws.Range("C32").Value = "ns"
$ws.Range("D32").Value = "ze"
$ws.Range("E32").Value = "ps"
$ws.Range("F32").Value = "pl"

# --- Rows 33-36: single-column values ---
$ws.Range("A33").Value = "a"
$ws.Range("A34").Value = "b"
$ws.Range("A35").Value = "c"
$ws.Range("A36").Value = "d"

# --- Update the saved selection to match the author's final cursor position ---
$ws.Range("L27").Select()
